# Update the cached "datetimeFigureOut" footer date (11/02/2019 -> 12/02/2019)
# on the slide master and on every slide layout, and remove the stray
# connector shape ("Conector recto de flecha 39", Id 40) from slide 1.

$p = $ppt.ActivePresentation

$oldDate = "11/02/2019"
$newDate = "12/02/2019"
$ppPlaceholderDate = 16

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every custom (slide) layout keeps its own cached copy of the date field.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# Individual slides can also carry their own cached copy of the date field
# (e.g. when the footer placeholder was overridden per-slide).
for ($sidx = 1; $sidx -le $p.Slides.Count; $sidx++) {
    Update-DateField $p.Slides.Item($sidx).Shapes
}

# Remove the orphan "Conector recto de flecha 39" connector (Id 40), wherever
# it shows up.
for ($sidx = 1; $sidx -le $p.Slides.Count; $sidx++) {
    $slide = $p.Slides.Item($sidx)
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq 40 -and $sh.Name -eq "Conector recto de flecha 39") {
            $sh.Delete()
        }
    }
}
